# CartaGantt.xlsx - progress update + view scroll
#
# Commit: "FIX: Se realizan validaciones de seguridad de la creación de
# contraseña. Se agregan paginas 404 500 y errores genéricos"
#
# The gantt sheet ("calendarioproyecto") tracks task completion (0..1) in
# column D, and a couple of HDU rows that ship the 404/500 error pages /
# generic-error handling are marked complete. The "Semana para mostrar"
# control cell (E4) is advanced by one week, and the window is left
# scrolled/zoomed to where the user was working (around row 44, column CR).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calendarioproyecto")

# --- Week-to-show control cell (named range SemanaParaMostrar) ---
$ws.Range("E4").Value = -7

# --- Task progress (column D, 0 = 0%, 1 = 100%) ---
# D30: "HDU E1 - Subida de evidencias (8 hrs)" 80% -> 100%
$ws.Range("D30").Value = 1
# D35: "HDU G1 - Responsividad móvil (4 hrs)" 0% -> 100%
$ws.Range("D35").Value = 1
# D39: "HDU H1 - Páginas de error 404/500 (4 hrs)" 0% -> 100%
$ws.Range("D39").Value = 1
# D40: "HDU H2 - Notificaciones/toasts (4 hrs)" 0% -> 100%
$ws.Range("D40").Value = 1

# --- Window view: zoom out and scroll the frozen pane toward the new work ---
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("CR44").Select()
